# "Generate Report for Handoff"
#
# The localization-status report was regenerated: the three files that were
# previously reported as "Handed back: in sync with en-US" are now
# "Ready for handoff" again, and the handoff timestamps were refreshed.
# Re-running the report also reflows the (auto-fitted) "Status" column to
# match the width of its new, shorter text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-05 13:13:48"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-05 13:13:44"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-05 13:13:48"

# --- Re-fit the "Status" columns now that their text is shorter ------------
# ("Ready for handoff" is much shorter than "Handed back: in sync with
# en-US", so the report's column autosizing shrinks these columns.)
$newStatusColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth
